# Update "想去人数" (interest count) figures for the new data pull.
# Matches commit: "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 21331
$ws1.Range("F6").Value  = 1132
$ws1.Range("F8").Value  = 7954
$ws1.Range("F14").Value = 189
$ws1.Range("F15").Value = 173
$ws1.Range("F19").Value = 1361
$ws1.Range("F20").Value = 540
$ws1.Range("F25").Value = 86
$ws1.Range("F28").Value = 55
$ws1.Range("F32").Value = 607
$ws1.Range("F34").Value = 142
$ws1.Range("F35").Value = 5071
$ws1.Range("F40").Value = 13152
$ws1.Range("F45").Value = 311
$ws1.Range("F46").Value = 442
$ws1.Range("F47").Value = 4066
$ws1.Range("F48").Value = 332

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 332

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 21331
$ws4.Range("F5").Value  = 1132
$ws4.Range("F7").Value  = 7954
$ws4.Range("F13").Value = 189
$ws4.Range("F14").Value = 173
$ws4.Range("F17").Value = 1361
$ws4.Range("F18").Value = 540
$ws4.Range("F23").Value = 86
$ws4.Range("F26").Value = 55
$ws4.Range("F29").Value = 332
$ws4.Range("F30").Value = 607
$ws4.Range("F33").Value = 142
$ws4.Range("F35").Value = 5071
$ws4.Range("F40").Value = 13152
$ws4.Range("F45").Value = 311
$ws4.Range("F46").Value = 442
$ws4.Range("F47").Value = 4066
$ws4.Range("F48").Value = 332
